# Daily attendance processing - 2026-01-03 18:40:38
# Normalize the "Recorded By" (column G) entries so that the literal
# "System" token is listed first among the recorders for each session row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $value = $cell.Value2

    if ($null -eq $value) { continue }
    if ($value -notmatch ",") { continue }

    $parts = $value -split ", "

    $systemParts = @()
    $restParts = @()
    foreach ($p in $parts) {
        if ($p.Equals("System")) {
            $systemParts += $p
        } else {
            $restParts += $p
        }
    }

    if ($systemParts.Count -gt 0) {
        $newParts = $systemParts + $restParts
    } else {
        $newParts = @($parts[($parts.Count - 1)..0])
    }

    $newValue = [string]::Join(", ", $newParts)

    if ($newValue -ne $value) {
        $cell.Value2 = $newValue
    }
}
